$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enemies")

# --- New row 12: "Spec-Ops" enemy entry ---
# Copy the formatting of the row above (row 11) down into row 12 first,
# then overwrite the values/formula so the new row keeps the table look.
$ws.Range("B11:H11").Copy()
$ws.Range("B12:H12").PasteSpecial(-4104)
$excel.CutCopyMode = $false

$ws.Range("B12").Value = "Spec-Ops"
$ws.Range("C12").Value = 35
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = "Pierce"
$ws.Range("H12").Formula = "=C12*((D12*40)/32)"

# --- Update existing values per the balance pass ---
$ws.Range("C5").Value = 100
$ws.Range("C8").Value = 180

# --- Move the active selection ---
$ws.Range("D15").Select()

$wb.Application.CalculateFull()
